$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the old column H (which becomes J):
#   H = percentage_bouwjaarklasse_tot_2000
#   I = percentage_bouwjaarklasse_vanaf_2000
$ws.Columns("H:I").Insert()

# New column headers
$ws.Cells.Item(1, 8).Value = "percentage_bouwjaarklasse_tot_2000"
$ws.Cells.Item(1, 9).Value = "percentage_bouwjaarklasse_vanaf_2000"

# Per-row data for the two new columns, in sheet row order (rows 2..93)
$hdata = @(
    ,@(74, 26)
    ,@(85, 15)
    ,@(100, 0)
    ,@(98, 2)
    ,@(83, 17)
    ,@(79, 21)
    ,@(96, 4)
    ,@(91, 9)
    ,@(96, 4)
    ,@(86, 14)
    ,@(91, 9)
    ,@(87, 13)
    ,@(98, 2)
    ,@(100, 0)
    ,@(60, 40)
    ,@(-99999999, -99999999)
    ,@(93, 7)
    ,@(99, 1)
    ,@(98, 2)
    ,@(66, 34)
    ,@(95, 5)
    ,@(96, 4)
    ,@(94, 6)
    ,@(99, 1)
    ,@(100, 0)
    ,@(-99999999, -99999999)
    ,@(88, 12)
    ,@(95, 5)
    ,@(100, 0)
    ,@(100, 0)
    ,@(28, 72)
    ,@(90, 10)
    ,@(-99999999, -99999999)
    ,@(88, 12)
    ,@(-99999999, -99999999)
    ,@(95, 5)
    ,@(91, 9)
    ,@(70, 30)
    ,@(98, 2)
    ,@(33, 67)
    ,@(2, 98)
    ,@(92, 8)
    ,@(88, 12)
    ,@(100, 0)
    ,@(99, 1)
    ,@(94, 6)
    ,@(100, 0)
    ,@(94, 6)
    ,@(-99999999, -99999999)
    ,@(94, 6)
    ,@(91, 9)
    ,@(95, 5)
    ,@(96, 4)
    ,@(98, 2)
    ,@(77, 23)
    ,@(92, 8)
    ,@(87, 13)
    ,@(10, 90)
    ,@(63, 37)
    ,@(79, 21)
    ,@(95, 5)
    ,@(89, 11)
    ,@(95, 5)
    ,@(94, 6)
    ,@(75, 25)
    ,@(100, 0)
    ,@(89, 11)
    ,@(94, 6)
    ,@(92, 8)
    ,@(99, 1)
    ,@(100, 0)
    ,@(97, 3)
    ,@(100, 0)
    ,@(50, 50)
    ,@(87, 13)
    ,@(79, 21)
    ,@(-99999999, -99999999)
    ,@(-99999999, -99999999)
    ,@(83, 17)
    ,@(-99999999, -99999999)
    ,@(66, 34)
    ,@(93, 7)
    ,@(90, 10)
    ,@(93, 7)
    ,@(-99999999, -99999999)
    ,@(81, 19)
    ,@(80, 20)
    ,@(85, 15)
    ,@(-99999999, -99999999)
    ,@(-99999999, -99999999)
    ,@(-99999999, -99999999)
    ,@(-99999999, -99999999)
)

$r = 2
foreach ($pair in $hdata) {
    $ws.Cells.Item($r, 8).Value = $pair[0]
    $ws.Cells.Item($r, 9).Value = $pair[1]
    $r = $r + 1
}

# Re-apply the AutoFilter so its stored range grows from A1:H93 to A1:J93
$ws.AutoFilterMode = $false
$ws.Range("A1:J93").AutoFilter()

# The hidden _FilterDatabase defined name also needs to track the new range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Rotterdam!_FilterDatabase") {
        $n.RefersTo = "=Rotterdam!`$A`$1:`$J`$93"
    }
}

# Approximate the new/target column widths for H, I and J
$ws.Columns("H").ColumnWidth = 30.498697916666668
$ws.Columns("I").ColumnWidth = 32.998697916666664
$ws.Columns("J").ColumnWidth = 21.830729166666668

# Match the saved selection/view state (columns H:I selected)
$ws.Range("H1:I1048576").Select()
